# Update the build timestamp embedded in the "version" string wherever it
# appears in the workbook (About sheet header/citation text, and the
# build_version column on the data sheet).

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
                $cell.Value = $val.Replace($oldStamp, $newStamp)
            }
        }
    }
}
